$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The block of 65 rows that is currently at A131:E195 (the most recent
# entries, added Spring 2019) needs to move up to directly follow the
# header row, i.e. become A2:E66. Everything that is currently at
# A2:E130 shifts down by 65 rows to become A67:E195. Rows 196 onward are
# untouched.

# 1. Remember the values of the block that is moving up.
$srcVal = $ws.Range("A131:E195").Value2

# 2. Open up 65 blank rows right after the header row. This pushes the
#    old A2:E130 block down to A67:E195 and also pushes the source block
#    (now stale/duplicate) down to A196:E260.
$ws.Range("A2:E66").EntireRow.Insert()

# 3. Write the remembered values into the newly freed rows.
$ws.Range("A2:E66").Value2 = $srcVal

# 4. Copy the date-column number formatting from a known-good date cell
#    onto the freshly written date column so it keeps the same style
#    (rather than Excel re-deriving a brand new number format style).
$ws.Range("D67").Copy()
$ws.Range("D2:D66").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 5. Delete the now-duplicated original copy of the block, which after
#    the insert in step 2 lives at rows 196:260.
$ws.Range("A196:E260").EntireRow.Delete()
